# Update the fit-statistics table (rows 3-6) on the active sheet with the
# recalculated values coming from the refreshed analysis/tables/graphs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (M2)
$ws.Range("B3").Value = -5318.503699629235
$ws.Range("C3").Value = 2612.133008721167
$ws.Range("D3").Value = 847.0612574218505
$ws.Range("F3").Value = 10699.00739925847
$ws.Range("G3").Value = 10852.48265708019

# Row 4 (M3)
$ws.Range("B4").Value = -5220.635254642758
$ws.Range("C4").Value = 2110.289734076519
$ws.Range("D4").Value = 731.4149608103575
$ws.Range("F4").Value = 10535.27050928552
$ws.Range("G4").Value = 10767.95880340232

# Row 5 (M4)
$ws.Range("B5").Value = -4968.576074253863
$ws.Range("C5").Value = 2102.087236060291
$ws.Range("D5").Value = 646.9182203152285
$ws.Range("F5").Value = 10063.15214850773
$ws.Range("G5").Value = 10375.05347891961

# Row 6 (M5)
$ws.Range("B6").Value = -4845.883879931207
$ws.Range("C6").Value = 1004.606216689852
$ws.Range("D6").Value = 565.4222894737685
$ws.Range("F6").Value = 9849.767759862414
$ws.Range("G6").Value = 10240.88212656938
$ws.Range("H6").Value = 0.1828341105110729
